$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The questionnaire was extended with a full new "Hogwarts house sorting"
# question/answer set (rows 3-16, columns A-E). The original single
# "In welk onderdeel..." question (with a typo) is replaced by the corrected
# wording further down the list, and every row now has its own question +
# 4 answers instead of only the trailing specialisation columns (F:I).

# Row 3
$ws.Range("A3").Value = 'In welk onderdeel van de magische wereld van Hogwarts ben jij het meest geïnteresseerd?'
$ws.Range("B3").Value = 'In duistere magie?'
$ws.Range("C3").Value = 'Het brouwen van diverse potions?'
$ws.Range("D3").Value = 'Het bedenken van spreuken?'
$ws.Range("E3").Value = 'Of wordt jij de coach van jouw zwerkbal team?'

# Row 4
$ws.Range("A4").Value = 'Welke positie speel jij in zwerkbal?'
$ws.Range("B4").Value = 'Drijvers(De Jagers dwarsbomen zodat zij kunnen scoren)'
$ws.Range("C4").Value = 'Jagers(Zoveel mogelijk punten scoren)'
$ws.Range("D4").Value = 'Zoeker(Op zoek gaan naar de golden snitch om de winst te bemachtigen)'
$ws.Range("E4").Value = 'Wachter(De drie ringen verdedigen van de jagers die willen scoren)'

# Row 5
$ws.Range("A5").Value = 'Hoe wil je gaan voor de winst?'
$ws.Range("B5").Value = 'Inspelen op de zwakke punten van de tegenstander?'
$ws.Range("C5").Value = 'Zoveel mogelijk doelpunten scoren om zo tijd te winnen voor de zoeker?'
$ws.Range("D5").Value = 'Zo snel mogelijk de snitch pakken?'
$ws.Range("E5").Value = 'Zorgen dat de topspelers van de tegenstander eruit worden gespeeld door de drijvers?'

# Row 6
$ws.Range("A6").Value = 'Waar zoek jij je avontuur in Hogwarts?'
$ws.Range("B6").Value = 'De mysteries opzoeken in het verboden bos?'
$ws.Range("C6").Value = 'De meest vreemde potions bedenken in de Potion classroom?'
$ws.Range("D6").Value = 'Te weten komen over spreueken en geschiedenis in de library?'
$ws.Range("E6").Value = 'Of begkijk je graag de zwerkbal wedstrijden?'

# Row 7
$ws.Range("A7").Value = 'Stel je krijgt een potion voorgeschoteld die jij moet na brouwen, hoe ga je hiermee aan de slag?'
$ws.Range("B7").Value = 'Je onderzoekt de potion doorgrondig om zo alle ingrediënten te achterhalen.'
$ws.Range("C7").Value = 'Je gooit random ingrediënten bij elkaar en hoopt dat je het in 1x goed krijgt.'
$ws.Range("D7").Value = 'Via trial en error proberen achter de ingrediënten te komen.'
$ws.Range("E7").Value = 'Je verzamelt van alle andere opties de uitkomst en gebruikt die informatie om de potion te maken.'

# Row 8
$ws.Range("A8").Value = 'De headmaster heeft je staf afgenomen, hoe krijg je je staf weer terug?'
$ws.Range("B8").Value = 'Kijken of je een back=up hebt, zo niet dan terughaalsoftware gebruiken of schijven restoren.'
$ws.Range("C8").Value = 'Bedenk een oplossing om je staf weer terug te pakken. (Bijvoorbeeld jezelf onzichtbaar maken en zo in het kantoor je staf terugpakken) '
$ws.Range("D8").Value = 'Wachten tot je je staf krijgt en in de tussentijd werken met een leenstaf'
$ws.Range("E8").Value = 'Analyseert de mogelijkeden hoe ze hem terug kunnen krijgen op de beste manier en gebruiken de best practice.'

# Row 9
$ws.Range("A9").Value = 'Stel je hebt een curse gekregen van de duistere magie, wat zou je doen om de curse van je af te halen? '
$ws.Range("B9").Value = 'Ik weet zelf hoe ik dit van mezelf moet afhalen, geef me maar deze ingrediënten en ik los het wel op. '
$ws.Range("C9").Value = 'Je ziet de curse in als een raadsel die je gaat proberen op te lossen.'
$ws.Range("D9").Value = 'Ga naar de leraar toe en vraag om hulp omdat je het niet weet.'
$ws.Range("E9").Value = 'In je boeken over duistere magie zoeken naar de curse die je hebt en hoe je ervan afkomt.'

# Row 10
$ws.Range("A10").Value = 'Een vriend van je is erachter gekomen hoe je onzichtbaar kan worden, jij wilt dit ook kunnen alleen die vriend verteld dit niet aan je, hoe los je dit op?'
$ws.Range("B10").Value = 'Alle bronnen van diegene onderzoeken naar informatie die kan leiden tot dit geheim.'
$ws.Range("C10").Value = 'Proberen op alternatieve manieren om onzichtbaar te worden.'
$ws.Range("D10").Value = 'Ga naar de library en zoek naar een boek over onzichtbaarheid.'
$ws.Range("E10").Value = 'Analyseer alle mogelijkheden die je vriend heeft kunnen gebruiken en probeer het op de manier te doen waarop je vriend het heeft kunnen doen. '

# Row 11
$ws.Range("A11").Value = 'Hoe kom je erachter wat de vloek des doods is?'
$ws.Range("B11").Value = 'Je hebt er in het verleden over gehoord en probeert contact op te zoeken met diegene die het je al heeft verteld.'
$ws.Range("C11").Value = 'Je maakt jezelf onzichtbaar om in de verboden afdeling van de library te gaan.'
$ws.Range("D11").Value = 'Je wilt dit niet eens weten, want het spreekt je niet aan.'
$ws.Range("E11").Value = 'Je gaat op zoek naar boeken in de library die er van alles over vertellen om een beeld te creëren wat het met je kan doen.'

# Row 12
$ws.Range("A12").Value = 'Je bent in de bibliotheek, naar wat voor boek zoek je?'
$ws.Range("B12").Value = '‘The mysteries of the forbidden forest volume I’ '
$ws.Range("C12").Value = '‘How to craft a broom 101’'
$ws.Range("D12").Value = '‘The art of creating spells volume I’'
$ws.Range("E12").Value = '‘The greatest quidditch matches of alltime volume I’'

# Row 13
$ws.Range("A13").Value = 'Stel je hebt vrije loop in het verboden bos, wat ga je doen?'
$ws.Range("B13").Value = 'Je gaat dieper het bos in om te vinden wat er zo verboden aan is.'
$ws.Range("C13").Value = 'Je zoekt inspiratie voor een duister verhaal dat je aan het schrijven bent.'
$ws.Range("D13").Value = 'Je zoekt je weg naar buiten want je bent bang.'
$ws.Range("E13").Value = 'Je gaat opzoek naar de monsters waarover je hebt gehoord.'

# Row 14
$ws.Range("A14").Value = 'Je komt per ongeluk voor de deur te staan waar je vanaf het begin van het jaar niet naar binnen mag, wat doe je?'
$ws.Range("B14").Value = 'Je opent simpelweg de deur en kijkt wat erachter staat.'
$ws.Range("C14").Value = 'Je zoekt naar een oplossing om van buitenaf naar binnen te kijken (met een spreuk of met een magisch voorwerp bijv.)'
$ws.Range("D14").Value = 'Je gaat door het sleutelgat kijken wat erachter zit, omdat je toch best nieuwsgierig bent.'
$ws.Range("E14").Value = 'Bekijkt de mogelijkheden die er zijn en brengt dit in kaart.'

# Row 15
$ws.Range("A15").Value = 'Je ziet iemand een spell uitvoeren, je ziet dat de kleur een blauwe flair heeft, hoe zou je achterhalen hoe de flair blauw is?'
$ws.Range("B15").Value = 'Je kan uitzoeken welke benodigdheden je nodig hebt om de de spel tot stand te krijgen, dit vervolgens noteren en rapport voor uitdraaien.'
$ws.Range("C15").Value = 'De kleur blauw spreekt je niet aan dus je gaat proberen spells te maken met een andere flair dan blauw.'
$ws.Range("D15").Value = 'Je zoekt die specifieke spell en achterhaalt daaruit wat de flair blauw maakt.'
$ws.Range("E15").Value = 'Je kijkt naar andere spells met een blauwe flair en probeer dit te reproduceren.'

# Row 16
$ws.Range("A16").Value = 'Je spreekt per ongeluk de naam uit van degene “Waarvan de naam niet mag worden genoemd uit”, wat doe je?'
$ws.Range("B16").Value = 'Je wacht af wat er gebeurt en je loopt rustig verder.'
$ws.Range("C16").Value = 'Kan kijken naar het gevolg en proberen hierop in te spelen'
$ws.Range("D16").Value = 'Ik bouw een muur om mij heen en isoleer mij van de buitenwereld.'
$ws.Range("E16").Value = 'Je geeft op omdat het de afgelopen keren ook in de dood is afgelopen.'

# A few answer cells were (re)typed with an explicit black font colour
# instead of the default automatic/theme colour.
$ws.Range("D10").Font.Color = 0
$ws.Range("E10").Font.Color = 0
$ws.Range("D11").Font.Color = 0

# Restore the page setup (paper size + orientation) and selection that were
# captured when the sheet was last saved.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("E28").Select()
